# Apply cryptocurrency price/volume updates per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '69.774.51'
$ws.Cells.Item(2, 5).Value = '  -1.50%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.498.45'
$ws.Cells.Item(3, 5).Value = '  -3.80%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.13%  '

# Row 5
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '581.06'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -4.18%  '

# Row 6
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '193.31'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -2.93%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '3.485.56'
$ws.Cells.Item(8, 5).Value = '  -3.77%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  -0.03%  '

# Row 10
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.205'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -6.93%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  -4.45%  '

# Row 12
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '51.71'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -4.37%  '

# Row 13
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000287'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -6.13%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  -4.38%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.031.37'
$ws.Cells.Item(15, 5).Value = '  -4.15%  '

# Row 16
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '647.43'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -5.97%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '69.568.32'
$ws.Cells.Item(17, 5).Value = '  -1.83%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.492.17'
$ws.Cells.Item(18, 5).Value = '  -3.72%  '

# Row 19
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '12.34'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -5.32%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.75%  '

# Row 21
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '18.24'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -4.10%  '

# Row 22
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '0.949'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -5.22%  '

# Row 23
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '18.29'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -2.99%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -2.56%  '

# Row 25
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '98.97'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -6.07%  '

# Row 26
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '4.28'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -7.44%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -4.00%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -3.95%  '

# Row 29
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '9.36'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -5.19%  '

# Row 30
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '32.68'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -4.86%  '

# Row 31
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '4.25'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -7.91%  '

# Row 32
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '6.75'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -6.22%  '

# Row 33
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = '11.65'
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -4.46%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -4.87%  '

# Row 35
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '61.37'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -3.17%  '

# Row 36
$c = $ws.Cells.Item(36, 4)
$c.NumberFormat = "@"
$c.Value = '531.82'
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +5.34%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '3.708.70'
$ws.Cells.Item(37, 5).Value = '  -6.23%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +0.11%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '0.0₃0788'
$ws.Cells.Item(39, 5).Value = '  -9.62%  '

# Row 40
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '3.56'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.15%  '

# Row 41
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '2.93'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -4.01%  '

# Row 42
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.374'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -3.82%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Kaspa'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '0.133'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -2.63%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '34.39'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -6.64%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'CoreDAO'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '3.50'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +31.22%  '

# Row 46
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.0443'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -3.45%  '

# Row 47
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -3.77%  '

# Row 48
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '2.83'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -8.41%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -4.26%  '

# Row 50
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.Style = "Normal"

# Row 51
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '8.20'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -5.58%  '
